{"js": "// Update the two-digit multiplication equations in the document body.\n// Each old equation text is unique in the document, so we can safely\n// find-and-replace each one independently.\nconst replacements = [\n  [\"66\u00d769=\", \"96\u00d757=\"],\n  [\"99\u00d736=\", \"38\u00d723=\"],\n  [\"26\u00d723=\", \"26\u00d740=\"],\n  [\"94\u00d712=\", \"28\u00d732=\"],\n  [\"39\u00d748=\", \"88\u00d715=\"],\n  [\"53\u00d753=\", \"92\u00d752=\"],\n  [\"54\u00d782=\", \"66\u00d724=\"],\n  [\"28\u00d722=\", \"85\u00d791=\"],\n  [\"46\u00d784=\", \"92\u00d779=\"],\n  [\"25\u00d787=\", \"68\u00d729=\"],\n  [\"12\u00d788=\", \"65\u00d758=\"],\n  [\"28\u00d786=\", \"46\u00d733=\"],\n  [\"30\u00d727=\", \"66\u00d725=\"],\n  [\"60\u00d757=\", \"74\u00d767=\"],\n  [\"32\u00d768=\", \"64\u00d752=\"],\n  [\"67\u00d779=\", \"15\u00d734=\"],\n  [\"38\u00d765=\", \"28\u00d792=\"],\n  [\"75\u00d724=\", \"68\u00d754=\"],\n  [\"76\u00d739=\", \"45\u00d739=\"],\n  [\"57\u00d741=\", \"89\u00d771=\"],\n  [\"83\u00d720=\", \"30\u00d745=\"],\n  [\"14\u00d737=\", \"27\u00d793=\"],\n  [\"82\u00d723=\", \"45\u00d719=\"],\n  [\"19\u00d722=\", \"34\u00d732=\"],\n  [\"19\u00d764=\", \"41\u00d731=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit multiplication equations in the document body.\n# Each old equation text is unique in the document, so a simple\n# Find/Replace (wdReplaceAll) for each pair is safe and deterministic.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"66\u00d769=\", \"96\u00d757=\"),\n    @(\"99\u00d736=\", \"38\u00d723=\"),\n    @(\"26\u00d723=\", \"26\u00d740=\"),\n    @(\"94\u00d712=\", \"28\u00d732=\"),\n    @(\"39\u00d748=\", \"88\u00d715=\"),\n    @(\"53\u00d753=\", \"92\u00d752=\"),\n    @(\"54\u00d782=\", \"66\u00d724=\"),\n    @(\"28\u00d722=\", \"85\u00d791=\"),\n    @(\"46\u00d784=\", \"92\u00d779=\"),\n    @(\"25\u00d787=\", \"68\u00d729=\"),\n    @(\"12\u00d788=\", \"65\u00d758=\"),\n    @(\"28\u00d786=\", \"46\u00d733=\"),\n    @(\"30\u00d727=\", \"66\u00d725=\"),\n    @(\"60\u00d757=\", \"74\u00d767=\"),\n    @(\"32\u00d768=\", \"64\u00d752=\"),\n    @(\"67\u00d779=\", \"15\u00d734=\"),\n    @(\"38\u00d765=\", \"28\u00d792=\"),\n    @(\"75\u00d724=\", \"68\u00d754=\"),\n    @(\"76\u00d739=\", \"45\u00d739=\"),\n    @(\"57\u00d741=\", \"89\u00d771=\"),\n    @(\"83\u00d720=\", \"30\u00d745=\"),\n    @(\"14\u00d737=\", \"27\u00d793=\"),\n    @(\"82\u00d723=\", \"45\u00d719=\"),\n    @(\"19\u00d722=\", \"34\u00d732=\"),\n    @(\"19\u00d764=\", \"41\u00d731=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
